$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output $ws.Name
Write-Output $ws.Range("AQ3").Formula
